# Generate Report for Handback
# Update the timestamp values that were refreshed when the handback report was regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for aa91b532-... row (also shared with de-de!H4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-17 02:48:56"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# "Correspond Handoff Datetime" for aa91b532-... row
$wsZhCn.Range("H4").Value = "2016-08-17 02:48:51"
# "Correspond Handback DateTime" for aa91b532-... row
$wsZhCn.Range("K4").Value = "2016-08-17 02:49:12"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# "Correspond Handback DateTime" for aa91b532-... row
$wsDeDe.Range("K4").Value = "2016-08-17 02:49:19"
